$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-06-19 Wednesday" "2024-06-20 Thursday"

Replace-Text "422×6=2532" "623×9=5607"
Replace-Text "153×6=918" "759×2=1518"
Replace-Text "353×4=1412" "593×2=1186"
Replace-Text "540×7=3780" "542×7=3794"
Replace-Text "695×6=4170" "772×9=6948"
Replace-Text "177×8=1416" "747×9=6723"
Replace-Text "159×2=318" "816×9=7344"
Replace-Text "286×3=858" "909×7=6363"
Replace-Text "391×9=3519" "121×8=968"
Replace-Text "735×2=1470" "810×3=2430"
Replace-Text "148×4=592" "902×6=5412"
Replace-Text "274×4=1096" "636×4=2544"
Replace-Text "957×6=5742" "940×3=2820"
Replace-Text "191×6=1146" "259×8=2072"
Replace-Text "606×3=1818" "231×3=693"
Replace-Text "403×8=3224" "304×2=608"
Replace-Text "189×4=756" "430×5=2150"
Replace-Text "701×8=5608" "779×7=5453"
Replace-Text "808×7=5656" "288×2=576"
Replace-Text "812×8=6496" "983×7=6881"
Replace-Text "348×2=696" "938×5=4690"
Replace-Text "963×3=2889" "883×3=2649"
Replace-Text "195×3=585" "494×2=988"
Replace-Text "337×4=1348" "316×8=2528"
Replace-Text "287×6=1722" "820×5=4100"

Write-Output "Done"
